# Processed Results - updated measures
# Adds "Std" (sample standard deviation) and "Relative std" (relative to the
# mean, expressed as a percentage) next to the existing summary statistics
# block (Min/Max, Q1/Median, Q3/IQR) on the Gyroscope sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New header labels (row 15), bold like the other summary headers (D1, D6/E6, D9/E9, D12/E12)
$ws.Range("D15").Value = "Std"
$ws.Range("D15").Font.Bold = $true
$ws.Range("E15").Value = "Relative std"
$ws.Range("E15").Font.Bold = $true

# New formulas (row 16): sample standard deviation of the data, and the
# standard deviation relative to the mean (E4), expressed as a percentage.
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "=(D16/E4)*100"

# Remove the duplicate / stale auto-generated chart-tracking defined names.
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# Reflect the last-used selection.
$ws.Range("E16").Select()
